$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.804.74'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '3.246.97'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.03%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.134'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.41%  '
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.415'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("D12").Value = '3.809.85'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.17%  '
$ws.Range("D15").Value = '67.793.80'
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").Value = '3.247.71'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  +2.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '379.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.31%  '
$ws.Range("E27").Value = '  +2.00%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.28%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +2.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.836'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("E41").Value = '  +6.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.59'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.89%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '345.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.59%  '
$ws.Range("E46").Value = '  +2.43%  '
$ws.Range("D47").Value = '2.619.60'
$ws.Range("E47").Value = '  -3.28%  '
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("E49").Value = '  -0.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.992'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.21%  '
